{"js": "// Replace the two split \"<id>...</id>\" run sequences with a single run\n// each, updating the enclosed identifier at the same time:\n//   <id>p089v_a1</id>  ->  <id>p089v_1</id>\n//   <id>p089v_a2</id>  ->  <id>p089v_2</id>\n// Word's Range.insertText(..., \"Replace\") collapses the matched range\n// into one run using the formatting of the range's leading run, which\n// matches the Courier-New/7f6000 styling already used for the <id>/</id>\n// tags in the document.\nconst replacements = [\n  { find: \"<id>p089v_a1</id>\", text: \"<id>p089v_1</id>\" },\n  { find: \"<id>p089v_a2</id>\", text: \"<id>p089v_2</id>\" }\n];\n\nfor (const { find, text } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(text, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Collapse the two split \"<id>...</id>\" run sequences into a single run\n# each while updating the enclosed identifier:\n#   <id>p089v_a1</id>  ->  <id>p089v_1</id>\n#   <id>p089v_a2</id>  ->  <id>p089v_2</id>\n# Word's Find/Replace rewrites the matched text as a single run using the\n# formatting of the first run in the match (Courier New / color 7f6000),\n# which mirrors the existing styling of the <id> and </id> tags.\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$found1 = $find1.Execute(\"<id>p089v_a1</id>\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"<id>p089v_1</id>\", $wdReplaceAll)\nWrite-Output \"Replaced p089v_a1: $found1\"\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$found2 = $find2.Execute(\"<id>p089v_a2</id>\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"<id>p089v_2</id>\", $wdReplaceAll)\nWrite-Output \"Replaced p089v_a2: $found2\"\n"}
